# Replacing the models for Kahraman with better ones
#
# - Shift every timestamp in column A (rows 2-97) forward by exactly one day.
# - Update the Notified Production values in column B for the rows whose
#   forecast changed (rows 23-73, 76-77).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps (column A, rows 2 through 97) forward by one day.
for ($r = 2; $r -le 97; $r++) {
    $cur = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value2 = $cur + 1
}

# New "Notified Production (MW)" values (column B) for the rows that changed.
$bUpdates = @{
    23 = 22
    24 = 25
    25 = 29
    26 = 204
    27 = 218
    28 = 237
    29 = 262
    30 = 542
    31 = 582
    32 = 616
    33 = 650
    34 = 920
    35 = 963
    36 = 1003
    37 = 1031
    38 = 1254
    39 = 1276
    40 = 1293
    41 = 1311
    42 = 1409
    43 = 1419
    44 = 1431
    45 = 1442
    46 = 1491
    47 = 1498
    48 = 1504
    49 = 1506
    50 = 1468
    51 = 1462
    52 = 1451
    53 = 1433
    54 = 1299
    55 = 1274
    56 = 1252
    57 = 1226
    58 = 1000
    59 = 968
    60 = 934
    61 = 898
    62 = 600
    63 = 563
    64 = 532
    65 = 506
    66 = 230
    67 = 212
    68 = 194
    69 = 181
    70 = 22
    71 = 16
    72 = 13
    73 = 12
    76 = 1
    77 = 1
}

foreach ($r in $bUpdates.Keys) {
    $ws.Cells.Item($r, 2).Value2 = $bUpdates[$r]
}
